$wb = $excel.ActiveWorkbook

# Delete Sheet2 - its data now fully lives in Sheet1 (dataprovider now drives
# filepath/sheetname dynamically, so a single sheet with more rows replaces
# the old two-sheet layout)
$wb.Worksheets("Sheet2").Delete()

$ws = $wb.Worksheets("Sheet1")

# Columns B-F first (these reuse existing shared strings plus the two new
# locale/timezone strings), then column A (the new venkat1..venkat4 ids) -
# this ordering reproduces the shared-string table layout of the target file.
$ws.Range("B1").Value = "p"
$ws.Range("C1").Value = "System - Administrator"
$ws.Range("D1").Value = "Administrator"
$ws.Range("E1").Value = "English (United States)"
$ws.Range("F1").Value = "(UTC+01:00:00) Europe/Amsterdam"

$ws.Range("B2").Value = "p"
$ws.Range("C2").Value = "System - Administrator"
$ws.Range("D2").Value = "Administrator"
$ws.Range("E2").Value = "Arabic (U.A.E.)"
$ws.Range("F2").Value = "(UTC+00:00:00) Atlantic/Canary"

$ws.Range("A1").Value = "venkat1"
$ws.Range("A2").Value = "venkat2"
$ws.Range("A3").Value = "venkat3"
$ws.Range("A4").Value = "venkat4"

$ws.Range("B3").Value = "p"
$ws.Range("C3").Value = "System - Administrator"
$ws.Range("D3").Value = "Administrator"
$ws.Range("E3").Value = "English (United States)"
$ws.Range("F3").Value = "(UTC+01:00:00) Europe/Amsterdam"

$ws.Range("B4").Value = "p"
$ws.Range("C4").Value = "System - Administrator"
$ws.Range("D4").Value = "Administrator"
$ws.Range("E4").Value = "Arabic (U.A.E.)"
$ws.Range("F4").Value = "(UTC+00:00:00) Atlantic/Canary"

# Move selection to A5, matching the post-edit cursor position
$ws.Range("A5").Select()
